$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update Maximo value in C2 ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("C2").Value = 627.6299599956752

# --- Sheet "Solucion": shuffle Pedido/Salida assignment rows 2-41 ---
$wsSolucion = $wb.Worksheets.Item("Solucion")

$wsSolucion.Range("A2").Value = "Pedido_9"
$wsSolucion.Range("B2").Value = "S001"
$wsSolucion.Range("A3").Value = "Pedido_2"
$wsSolucion.Range("B3").Value = "S025"
$wsSolucion.Range("A4").Value = "Pedido_4"
$wsSolucion.Range("B4").Value = "S005"
$wsSolucion.Range("A5").Value = "Pedido_18"
$wsSolucion.Range("B5").Value = "S029"
$wsSolucion.Range("A6").Value = "Pedido_20"
$wsSolucion.Range("B6").Value = "S002"
$wsSolucion.Range("A7").Value = "Pedido_5"
$wsSolucion.Range("B7").Value = "S026"
$wsSolucion.Range("A8").Value = "Pedido_30"
$wsSolucion.Range("B8").Value = "S006"
$wsSolucion.Range("A9").Value = "Pedido_37"
$wsSolucion.Range("B9").Value = "S003"
$wsSolucion.Range("A10").Value = "Pedido_22"
$wsSolucion.Range("B10").Value = "S030"
$wsSolucion.Range("A11").Value = "Pedido_7"
$wsSolucion.Range("B11").Value = "S007"
$wsSolucion.Range("A12").Value = "Pedido_38"
$wsSolucion.Range("B12").Value = "S027"
$wsSolucion.Range("A13").Value = "Pedido_16"
$wsSolucion.Range("B13").Value = "S031"
$wsSolucion.Range("A14").Value = "Pedido_28"
$wsSolucion.Range("B14").Value = "S004"
$wsSolucion.Range("A15").Value = "Pedido_23"
$wsSolucion.Range("B15").Value = "S008"
$wsSolucion.Range("A16").Value = "Pedido_13"
$wsSolucion.Range("B16").Value = "S028"
$wsSolucion.Range("A17").Value = "Pedido_17"
$wsSolucion.Range("B17").Value = "S009"
$wsSolucion.Range("A18").Value = "Pedido_11"
$wsSolucion.Range("B18").Value = "S032"
$wsSolucion.Range("A19").Value = "Pedido_14"
$wsSolucion.Range("B19").Value = "S013"
$wsSolucion.Range("A20").Value = "Pedido_3"
$wsSolucion.Range("B20").Value = "S033"
$wsSolucion.Range("A21").Value = "Pedido_19"
$wsSolucion.Range("B21").Value = "S010"
$wsSolucion.Range("A22").Value = "Pedido_6"
$wsSolucion.Range("B22").Value = "S037"
$wsSolucion.Range("A23").Value = "Pedido_40"
$wsSolucion.Range("B23").Value = "S014"
$wsSolucion.Range("A24").Value = "Pedido_35"
$wsSolucion.Range("B24").Value = "S034"
$wsSolucion.Range("A25").Value = "Pedido_21"
$wsSolucion.Range("B25").Value = "S011"
$wsSolucion.Range("A26").Value = "Pedido_29"
$wsSolucion.Range("B26").Value = "S015"
$wsSolucion.Range("A27").Value = "Pedido_25"
$wsSolucion.Range("B27").Value = "S038"
$wsSolucion.Range("A28").Value = "Pedido_26"
$wsSolucion.Range("B28").Value = "S012"
$wsSolucion.Range("A29").Value = "Pedido_36"
$wsSolucion.Range("B29").Value = "S016"
$wsSolucion.Range("A30").Value = "Pedido_12"
$wsSolucion.Range("B30").Value = "S035"
$wsSolucion.Range("A31").Value = "Pedido_34"
$wsSolucion.Range("B31").Value = "S017"
$wsSolucion.Range("A32").Value = "Pedido_24"
$wsSolucion.Range("B32").Value = "S039"
$wsSolucion.Range("A33").Value = "Pedido_1"
$wsSolucion.Range("B33").Value = "S021"
$wsSolucion.Range("A34").Value = "Pedido_33"
$wsSolucion.Range("B34").Value = "S036"
$wsSolucion.Range("A35").Value = "Pedido_10"
$wsSolucion.Range("B35").Value = "S018"
$wsSolucion.Range("A36").Value = "Pedido_27"
$wsSolucion.Range("B36").Value = "S040"
$wsSolucion.Range("A37").Value = "Pedido_15"
$wsSolucion.Range("B37").Value = "S022"
$wsSolucion.Range("A38").Value = "Pedido_31"
$wsSolucion.Range("B38").Value = "S019"
$wsSolucion.Range("A39").Value = "Pedido_8"
$wsSolucion.Range("B39").Value = "S023"
$wsSolucion.Range("A40").Value = "Pedido_39"
$wsSolucion.Range("B40").Value = "S020"
$wsSolucion.Range("A41").Value = "Pedido_32"
$wsSolucion.Range("B41").Value = "S024"

# --- Sheet "Metricas": update Tiempo values for Z1 and Z2 ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 627.6299599956752
$wsMetricas.Range("B3").Value = 489.4081252027246
